$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet "Sheet1" -> "General" ---
$ws.Name = "General"

# --- Update description text in C5 (add "paramétrages des") ---
$ws.Range("C5").Value = "Prise en compte des vues 3D pour les objets traditionnels de construction (paramétrages des murs, toits, dalles, etc.)"

# --- Update description text in C6 (style Standard -> modification du style Standard présent par défaut) ---
$ws.Range("C6").Value = "Ajout des hachures 2D (vue de dessus) pour les masses élémentaires (modification du style Standard présent par défaut)"

# --- Add new row 7: version name (B7) + long description (C7) ---
# Clone formatting from row 6 (same border / base alignment), then tweak alignment per cell.
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = "_C3D-TEMPLATE_2025_FRA (Architecture v0001c)"
$ws.Range("B7").VerticalAlignment = -4108

$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$longText = @"
Ajout d'un style Standard ACA pour les objets architecturaux suivants : 
- Définitions de formes d'éléments de structure ;
- Définitions des groupes de nettoyage de murs ;
- Styles de blocs porte/fenêtre ;
- Styles de bords de dalle ;
- Styles de bords de dalle du toit ;
- Styles de dalles ;
- Styles de dalles du toit ;
- Styles de fenêtres ;
- Styles de garde-corps ;
- Styles de murs ;
- Styles de murs-rideaux ;
- Styles de portes ;
- Styles d'escaliers ;
- Styles d'espaces (avec en supplément les styles Commerce, Education, Habitation, Logement, Santé) ;
- Styles d'extrémités de mur ;
- Styles d'extrémités d'ouverture de mur ;
- Styles d'unités de mur-rideau.
L'implémentation de ces styles Standard ACA dans le gabarit empêche Civil 3D de créer ses propres styles Standard par défaut qui seraient mal paramétrés ou avec des composants manquants dans l'affichage (les rambardes des garde-corps par exemple).
Les styles de balancements d'escaliers Equilibré, Manuel, Sur un point sont également implémentés dans cette nouvelle version du gabarit.
"@
$ws.Range("C7").Value = $longText
$ws.Range("C7").WrapText = $true

# --- Row height for the new (tall, wrapped) row ---
$ws.Rows(7).RowHeight = 299.25

# --- Widen column C to fit the long description ---
$ws.Columns("C").ColumnWidth = 112.5

# --- Selection as left by the author ---
$ws.Range("C13").Select() | Out-Null
